$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6189727187156677
$ws.Range("B1").Value = 0.8037946820259094
$ws.Range("C1").Value = 1.237553477287292
$ws.Range("D1").Value = 5.565415382385254
$ws.Range("E1").Value = 4.533619403839111
